$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row / test case id
$ws.Range("A6").Value = "createnewpaymentrequest_ID"

# New header cells for payment request columns
$ws.Range("G1").Value = "Title"
$ws.Range("H1").Value = "Description"

$ws.Range("G6").Value = "Payment 1"
$ws.Range("H6").Value = "This is Test"

$ws.Range("I1").Value = "Amount"
$ws.Range("I6").Value = 500

$ws.Range("G1:I1").Interior.Color = $ws.Range("A1").Interior.Color

$ws.Columns.Item(1).ColumnWidth = 28.6

$ws.Range("I6").Select()
